# Apply weekly price-list rotation for Fruta/Pomelo Lo Valledor sheet
# (commit: "Fruta / hortaliza, semanal") — each data row (2..23) is updated
# to the values that, per the canonical diff, now belong to it.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44400
$ws.Range("M2").Value = 140
$ws.Range("N2").Value = 9800
$ws.Range("O2").Value = 9800
$ws.Range("P2").Value = 9800
$ws.Range("Q2").Value = '$/caja 14 kilos empedrada'
$ws.Range("S2").Value = 700
$ws.Range("T2").Value = 14

# Row 3
$ws.Range("D3").Value = 44193
$ws.Range("M3").Value = 8
$ws.Range("N3").Value = 150000
$ws.Range("O3").Value = 150000
$ws.Range("P3").Value = 150000
$ws.Range("S3").Value = 429

# Row 4
$ws.Range("D4").Value = 44356
$ws.Range("L4").Value = 'Primera'
$ws.Range("M4").Value = 24
$ws.Range("O4").Value = 230000
$ws.Range("P4").Value = 215000
$ws.Range("R4").Value = 'Región Metropolitana'
$ws.Range("S4").Value = 614

# Row 5
$ws.Range("D5").Value = 44189
$ws.Range("M5").Value = 16
$ws.Range("N5").Value = 150000
$ws.Range("P5").Value = 150000
$ws.Range("R5").Value = 'Provincia de Limarí'
$ws.Range("S5").Value = 429

# Row 6
$ws.Range("D6").Value = 44363
$ws.Range("M6").Value = 20
$ws.Range("N6").Value = 200000
$ws.Range("O6").Value = 230000
$ws.Range("P6").Value = 215000
$ws.Range("R6").Value = 'Provincia de Limarí'
$ws.Range("S6").Value = 614

# Row 7
$ws.Range("D7").Value = 44201
$ws.Range("L7").Value = 'Especial'
$ws.Range("M7").Value = 8
$ws.Range("N7").Value = 200000
$ws.Range("O7").Value = 200000
$ws.Range("P7").Value = 200000
$ws.Range("R7").Value = 'Región de O''Higgins'
$ws.Range("S7").Value = 571

# Row 8
$ws.Range("D8").Value = 44201
$ws.Range("N8").Value = 170000
$ws.Range("O8").Value = 170000
$ws.Range("P8").Value = 170000
$ws.Range("R8").Value = 'Región de O''Higgins'
$ws.Range("S8").Value = 486

# Row 9
$ws.Range("D9").Value = 44505
$ws.Range("M9").Value = 15
$ws.Range("N9").Value = 150000
$ws.Range("O9").Value = 150000
$ws.Range("P9").Value = 150000
$ws.Range("R9").Value = 'Provincia de Quillota'
$ws.Range("S9").Value = 429

# Row 10
$ws.Range("D10").Value = 44309
$ws.Range("L10").Value = 'Primera'
$ws.Range("N10").Value = 350000
$ws.Range("O10").Value = 350000
$ws.Range("P10").Value = 350000
$ws.Range("R10").Value = 'Región Metropolitana'
$ws.Range("S10").Value = 1000

# Row 11
$ws.Range("D11").Value = 44312
$ws.Range("L11").Value = 'Segunda'
$ws.Range("M11").Value = 10
$ws.Range("N11").Value = 330000
$ws.Range("O11").Value = 330000
$ws.Range("P11").Value = 330000
$ws.Range("R11").Value = 'Región Metropolitana'
$ws.Range("S11").Value = 943

# Row 12
$ws.Range("D12").Value = 44376
$ws.Range("N12").Value = 180000
$ws.Range("O12").Value = 180000
$ws.Range("P12").Value = 180000
$ws.Range("R12").Value = 'Hijuelas'
$ws.Range("S12").Value = 514

# Row 13
$ws.Range("D13").Value = 44376
$ws.Range("L13").Value = 'Segunda'
$ws.Range("M13").Value = 16
$ws.Range("N13").Value = 150000
$ws.Range("O13").Value = 150000
$ws.Range("P13").Value = 150000
$ws.Range("Q13").Value = '$/bins (350 kilos)'
$ws.Range("R13").Value = 'Provincia de Limarí'
$ws.Range("S13").Value = 429
$ws.Range("T13").Value = 350

# Row 14
$ws.Range("D14").Value = 44167
$ws.Range("M14").Value = 140
$ws.Range("N14").Value = 9800
$ws.Range("O14").Value = 9800
$ws.Range("P14").Value = 9800
$ws.Range("Q14").Value = '$/caja 14 kilos empedrada'
$ws.Range("R14").Value = 'Región de O''Higgins'
$ws.Range("S14").Value = 700
$ws.Range("T14").Value = 14

# Row 15
$ws.Range("D15").Value = 44208
$ws.Range("K15").Value = 'Start Ruby'
$ws.Range("M15").Value = 16
$ws.Range("N15").Value = 180000
$ws.Range("O15").Value = 180000
$ws.Range("P15").Value = 180000
$ws.Range("R15").Value = 'Región Metropolitana'
$ws.Range("S15").Value = 514

# Row 16
$ws.Range("D16").Value = 44389
$ws.Range("L16").Value = 'Especial'
$ws.Range("M16").Value = 18
$ws.Range("N16").Value = 200000
$ws.Range("O16").Value = 200000
$ws.Range("P16").Value = 200000
$ws.Range("S16").Value = 571

# Row 17
$ws.Range("D17").Value = 44308
$ws.Range("L17").Value = 'Primera'
$ws.Range("M17").Value = 20
$ws.Range("N17").Value = 280000
$ws.Range("O17").Value = 280000
$ws.Range("P17").Value = 280000
$ws.Range("R17").Value = 'Región Metropolitana'
$ws.Range("S17").Value = 800

# Row 18
$ws.Range("D18").Value = 44446
$ws.Range("M18").Value = 14
$ws.Range("N18").Value = 150000
$ws.Range("O18").Value = 160000
$ws.Range("P18").Value = 155000
$ws.Range("S18").Value = 443

# Row 19
$ws.Range("D19").Value = 44196
$ws.Range("K19").Value = 'Red Blush'
$ws.Range("M19").Value = 12
$ws.Range("N19").Value = 130000
$ws.Range("O19").Value = 130000
$ws.Range("P19").Value = 130000
$ws.Range("R19").Value = 'Provincia de Limarí'
$ws.Range("S19").Value = 371

# Row 22
$ws.Range("D22").Value = 44195
$ws.Range("M22").Value = 20
$ws.Range("N22").Value = 200000
$ws.Range("O22").Value = 210000
$ws.Range("P22").Value = 206000
$ws.Range("Q22").Value = '$/bins (350 kilos)'
$ws.Range("S22").Value = 589
$ws.Range("T22").Value = 350

# Row 23
$ws.Range("D23").Value = 44511
$ws.Range("M23").Value = 24
$ws.Range("N23").Value = 140000
$ws.Range("P23").Value = 145000
$ws.Range("R23").Value = 'Región Metropolitana'
$ws.Range("S23").Value = 414
